# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" positioned right before "总计"
#    (i.e. right after "2021-Q4"), populated with the quarterly fund
#    holdings detail.
# 2. Update the "总计" (totals) sheet with a new leading row for
#    2022-Q1 and shift the previously-existing rows down by one.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, $val) {
    # Forces a numeric-looking string ("43.52", "001882", ...) to be
    # stored as literal text instead of being auto-coerced to a number
    # (which would also destroy meaningful leading zeros in fund codes).
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$donorSheet = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Reuse the existing header/column-A cell formatting (bold + border,
# style index already present in the workbook) instead of creating new
# style entries. Copy header row and index column separately so no
# spurious empty A1 cell gets created (the source sheets never have one).
$donorSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

$donorSheet.Range("A2:A17").Copy()
$newSheet.Range("A2:A17").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holding rows (A=index, B=code, C=name, D=size, E=stock pos,
# F=pos pct, G=held value, H=pos rank)
$rows = @(
    @(0,  "166019", "中欧价值智选回报混合A",           "156.17", "94.14", "2.91", "4.5445", 10),
    @(1,  "166005", "中欧价值发现混合 -A",             "43.52",  "93.97", "7.44", "3.2379", 2),
    @(2,  "001882", "中欧价值发现混合 -E",             "43.52",  "93.97", "7.44", "3.2379", 2),
    @(3,  "001810", "中欧潜力价值灵活配置混合A",        "28.67",  "94.05", "7.75", "2.2219", 1),
    @(4,  "004235", "中欧价值智选回报混合C",           "36.40",  "94.14", "2.91", "1.0592", 10),
    @(5,  "004232", "中欧价值发现混合 -C",             "10.98",  "93.97", "7.44", "0.8169", 2),
    @(6,  "001887", "中欧价值智选回报混合E",           "20.77",  "94.14", "2.91", "0.6044", 10),
    @(7,  "004848", "中欧睿泓定期开放灵活配置混合",      "23.30",  "59.08", "2.06", "0.4800", 9),
    @(8,  "166024", "中欧恒利三年定期开放混合",         "4.48",   "98.71", "6.86", "0.3073", 2),
    @(9,  "005764", "中欧潜力价值灵活配置混合C",        "3.43",   "94.05", "7.75", "0.2658", 1),
    @(10, "001891", "中欧成长优选回报灵活配置混合E",     "2.97",   "94.42", "7.26", "0.2156", 1),
    @(11, "166020", "中欧成长优选回报灵活配置混合A",     "2.97",   "94.42", "7.26", "0.2156", 1),
    @(12, "004756", "国寿安保稳吉混合A",               "6.23",   "25.19", "0.78", "0.0486", 10),
    @(13, "710301", "富安达增强收益债券A",             "0.61",   "20.20", "4.67", "0.0285", 1),
    @(14, "004757", "国寿安保稳吉混合C",               "1.89",   "25.19", "0.78", "0.0147", 10),
    @(15, "710302", "富安达增强收益债券C",             "0.26",   "20.20", "4.67", "0.0121", 1)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $newSheet.Cells.Item($r, 1).Value = $row[0]          # A: index (number, inherits style 2 from the format copy)

    Set-TextCell $newSheet.Cells.Item($r, 2) $row[1]     # B: fund code (text)
    Set-TextCell $newSheet.Cells.Item($r, 3) $row[2]     # C: fund name (text)
    Set-TextCell $newSheet.Cells.Item($r, 4) $row[3]     # D: fund size (text)
    Set-TextCell $newSheet.Cells.Item($r, 5) $row[4]     # E: stock position (text)
    Set-TextCell $newSheet.Cells.Item($r, 6) $row[5]     # F: position pct (text)
    Set-TextCell $newSheet.Cells.Item($r, 7) $row[6]     # G: held value (text)

    $newSheet.Cells.Item($r, 8).Value = $row[7]          # H: position rank (number)
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: add a new leading row for 2022-Q1 and
#    shift the existing rows down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Extend the existing column-A style (bold + border) down into the new
# row 7 before overwriting values, so every data row A2:A7 keeps the
# same formatting the original A2:A6 rows had.
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)  # xlPasteFormats

$totalRows = @(
    @(0, "2022-Q1", 16, 17.31),
    @(1, "2021-Q4", 13, 5.56),
    @(2, "2021-Q3", 32, 14.16),
    @(3, "2021-Q2", 27, 10.91),
    @(4, "2021-Q1", 23, 6.86),
    @(5, "2020-Q4", 8, 6.19)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]

    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
}

Write-Output "2022-Q1 sheet added and 总计 sheet updated"
